# V3 Gerbers and STEP Files
# Rename the BOM shared-string title (V10 -> D16 revision) shown in the
# merged title cell D1:F3, and move the selection to that cell/range -
# matching the author re-selecting/editing the title after renaming the
# folder from .PCBsV2 to .PCBsV3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "T-1D16W547848A BOM"
$ws.Range("D1:F3").Select() | Out-Null
